$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "mejores arcades" -> "mejores arcade" (wrapped with gramStart/gramEnd
#    proofErr markers, splitting the run in three as Word's grammar checker
#    would do while leaving the rest of the sentence untouched).
# ---------------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("mejores arcades")
if (-not $found) {
  throw "Could not find 'mejores arcades'"
}

$para = $findRng.Paragraphs(1)
$pRng = $para.Range

$pkgNs = 'xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"'
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$objetivoXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="23215B2B" w14:textId="79D6A26A" w:rsidR="009A6FE0" w:rsidRDefault="000E5FF0" w:rsidP="000E5FF0" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:ind w:left="720"/><w:rPr><w:lang w:val="es-DO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-DO"/></w:rPr><w:t xml:space="preserve">El objetivo general de este proyecto es volver a recrear uno de los </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-DO"/></w:rPr><w:t>mejores arcade</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-DO"/></w:rPr><w:t xml:space="preserve"> de todos los tiempos y colocarlo en un espacio moderno donde los fan&#225;ticos de este tipo de video juegos puedan descargarlo y jugarlo en sus tiempos de entretenimiento.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$pRng.InsertXML($objetivoXml)

# ---------------------------------------------------------------------------
# 2) Add the five "Objetivos especificos" bullet items right after the
#    "Objetivos especificos" heading (before the existing blank paragraph
#    that precedes "Escenario").
# ---------------------------------------------------------------------------
$headRng = $d.Content
$headFound = $headRng.Find.Execute("Objetivos específicos")
if (-not $headFound) {
  throw "Could not find 'Objetivos específicos' heading"
}

$headPara = $headRng.Paragraphs(1)
$insertRng = $headPara.Range
$insertRng.Collapse(0)

$items = @(
  "Desarrollar o tratar de crear reflejos en los jugadores.",
  "Alcanzar la mayor cantidad de puntos y con esto competir con tus amigos.",
  "Entretenimiento simple y sencillo",
  "A pesar de ser un juego simple, dar un poco de reto a los jugadores.",
  "Dar un poco de nostalgia a los jugadores veteranos."
)

$bodyXml = ""
foreach ($item in $items) {
  $bodyXml += '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="48"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="8647"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Helvetica"/><w:szCs w:val="24"/><w:lang w:val="es-DO"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Helvetica"/><w:szCs w:val="24"/><w:lang w:val="es-DO"/></w:rPr><w:t>' + $item + '</w:t></w:r></w:p>'
}

$listXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRng.InsertXML($listXml)
